$d = $word.ActiveDocument

# --- Clear all existing content, leaving a single empty paragraph ---
$full = $d.Range(0, $d.Content.End)
$full.Delete()

# --- Pre-create all 7 empty paragraphs first, so that later formatting
#     (e.g. Bold on paragraph 1) does not leak/inherit into paragraphs
#     that are created afterwards via InsertParagraphAfter. ---
$p = $d.Paragraphs.Item(1)
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Item(2)
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Item(3)
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Item(4)
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Item(5)
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Item(6)
$p.Range.InsertParagraphAfter()

Write-Host "Paragraphs after pre-create: " $d.Paragraphs.Count

# --- Helper to write one centered paragraph of body text ---
function Set-CenteredParagraph($para, [string]$text, [bool]$bold, [int]$size) {
    if ($text -ne $null) {
        $para.Range.Text = $text
    }
    $para.Alignment = 1
    if ($bold) {
        $para.Range.Font.Bold = $true
        $para.Range.Font.BoldBi = $true
    }
    $para.Range.Font.Size = $size
    $para.Range.Font.SizeBi = $size
}

# 1) CERTIFICADO - big bold title
Set-CenteredParagraph $d.Paragraphs.Item(1) "CERTIFICADO" $true 36

# 2) Certificamos para os devidos fins que o colaborador:
Set-CenteredParagraph $d.Paragraphs.Item(2) "Certificamos para os devidos fins que o colaborador:" $false 16

# 3) Beatriz Costa
Set-CenteredParagraph $d.Paragraphs.Item(3) "Beatriz Costa" $false 16

# 4) Portador do CPF nº 999.888.777-66
Set-CenteredParagraph $d.Paragraphs.Item(4) "Portador do CPF nº 999.888.777-66" $false 16

# 5) Concluiu com êxito o treinamento de Norma Regulamentadora NR06.
Set-CenteredParagraph $d.Paragraphs.Item(5) "Concluiu com êxito o treinamento de Norma Regulamentadora NR06." $false 16

# 6) Rio de Janeiro, 15/01/2026.
Set-CenteredParagraph $d.Paragraphs.Item(6) "Rio de Janeiro, 15/01/2026." $false 16

# 7) Trailing empty (bold-formatted) paragraph
Set-CenteredParagraph $d.Paragraphs.Item(7) $null $true 16

# --- Add a single-line page border around the page (top/left/bottom/right) ---
$borders = $d.Sections.Item(1).Borders
$topBorder = $borders.Item(-1)
$topBorder.LineStyle = 1
$topBorder.LineWidth = 2

Write-Host "Paragraphs count: " $d.Paragraphs.Count
